# Add new translation entries for the "bouton deco" / "bouton admin" login
# buttons: login-button_connecte / Connecté, bouton_admin / Page d'administration,
# bouton_deco / Se déconnecter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$lo = $ws.ListObjects.Item("Traductions")

# Grow the "Traductions" table by three rows (A1:B51 -> A1:B54), as Excel
# does automatically when new rows are appended to a table.
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Fill the new rows. The order of entry below matches the order the new
# strings were appended to the shared-string table in the saved workbook:
# A52, A53, A54, B53, B54, B52.
$ws.Range("A52").Value = "login-button_connecte"
$ws.Range("A53").Value = "bouton_admin"
$ws.Range("A54").Value = "bouton_deco"

$ws.Range("B53").Value = "Page d'administration"
$ws.Range("B54").Value = "Se déconnecter"
$ws.Range("B52").Value = "Connecté"

# Restore the scroll position / selection left by the editing session.
$excel.ActiveWindow.ScrollRow = 35
$ws.Range("B53").Select()

$wb.Save()
